# daily auto push: 2026-02-27 07:10 UTC
# A new sample row for 2026/02/27 (weekday 金, hour 14) was logged. It sits
# right after the existing 2026/02/27 rows (885..888 originally) and pushes
# every subsequent row down by one (885 -> 886, ..., 926 -> 927).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: insert a blank row at 885, shifting
# 885:926 down to 886:927 (dimension grows from D926 to D927 automatically).
$ws.Rows("885:885").Insert()

# Column A stores dates as plain text (e.g. "2026/02/27"), not real Excel
# dates. Typing that string straight into Value would get auto-parsed into
# a date serial, so instead copy the already-text value from the row right
# above (A884, which holds the same "2026/02/27" text) using values-only
# paste - this keeps it a plain text cell with no extra formatting.
$ws.Range("A884").Copy()
$ws.Range("A885").PasteSpecial(-4163)   # xlPasteValues

# Remaining columns for the new row.
$ws.Range("B885").Value = "金"
$ws.Range("C885").Value = 14
$ws.Range("D885").Value = 201
